# Fruta / hortaliza, semanal
# Insert a new daily-price record row at row 342 (new weekly observation),
# pushing the existing rows 342:366 down to 343:367.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 342; this shifts rows 342:366 -> 343:367
$ws.Rows("342:342").Insert()

# Populate the new row 342. Columns A,B,C,E,F,G,H,I,K,L,M,N,O,P,Q,R are
# copied from the row immediately below (which now holds the data that used
# to be in row 342), while the date (D) and volume (J) get their new values.
$ws.Range("A342").Value2 = $ws.Range("A343").Value2
$ws.Range("B342").Value2 = $ws.Range("B343").Value2
$ws.Range("C342").Value2 = $ws.Range("C343").Value2
$ws.Range("D342").Value2 = 44746
$ws.Range("E342").Value2 = $ws.Range("E343").Value2
$ws.Range("F342").Value2 = $ws.Range("F343").Value2
$ws.Range("G342").Value2 = $ws.Range("G343").Value2
$ws.Range("H342").Value2 = $ws.Range("H343").Value2
$ws.Range("I342").Value2 = $ws.Range("I343").Value2
$ws.Range("J342").Value2 = 65
$ws.Range("K342").Value2 = $ws.Range("K343").Value2
$ws.Range("L342").Value2 = $ws.Range("L343").Value2
$ws.Range("M342").Value2 = $ws.Range("M343").Value2
$ws.Range("N342").Value2 = $ws.Range("N343").Value2
$ws.Range("O342").Value2 = $ws.Range("O343").Value2
$ws.Range("P342").Value2 = $ws.Range("P343").Value2
$ws.Range("Q342").Value2 = $ws.Range("Q343").Value2
$ws.Range("R342").Value2 = $ws.Range("R343").Value2
